$wb = $excel.ActiveWorkbook

# The edit happens on the "tryEditorCode" worksheet: cell A3's text is replaced
# with a new (truncated) snippet, and that sheet becomes the active/selected tab
# with A3 as the active cell.
$ws = $wb.Worksheets.Item("tryEditorCode")
$ws.Activate()

$ws.Range("A3").Value = "`tprint('Hello"

$ws.Range("A3").Select()
